$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (plain decimals like "19.15"); force them to remain plain text first,
# matching the inline-string / text cells used throughout this sheet.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

$ws.Range('D2').Value = '26.740.64'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '1.641.18'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').Value = '217.85'
$ws.Range('E5').Value = '  +1.57%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('D10').Value = '19.15'
$ws.Range('E10').Value = '  +0.34%  '
$ws.Range('D11').Value = '0.0843'
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').Value = '1.868.78'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '1.645.93'
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('E15').Value = '  -0.20%  '
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').Value = '26.733.28'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('D19').Value = '215.05'
$ws.Range('E19').Value = '  -0.10%  '
$ws.Range('D21').Value = '4.38'
$ws.Range('E21').Value = '  +1.07%  '
$ws.Range('D22').Value = '2.37'
$ws.Range('E22').Value = '  +7.06%  '
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('E24').Value = '  -1.71%  '
$ws.Range('D25').Value = '145.50'
$ws.Range('E25').Value = '  +0.36%  '
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('D27').Value = '0.119'
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('E28').Value = '  +0.82%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('D30').Value = '0.0509'
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('E31').Value = '  +1.58%  '
$ws.Range('E32').Value = '  +0.96%  '
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('D34').Value = '1.288.60'
$ws.Range('E34').Value = '  +0.81%  '
$ws.Range('D35').Value = '1.54'
$ws.Range('E35').Value = '  +0.31%  '
$ws.Range('E36').Value = '  +1.21%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('E38').Value = '  +1.17%  '
$ws.Range('D39').Value = '0.818'
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('E43').Value = '  -2.20%  '
$ws.Range('D44').Value = '1.779.30'
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').Value = '61.11'
$ws.Range('E45').Value = '  +3.15%  '
$ws.Range('D46').Value = '91.80'
$ws.Range('E46').Value = '  +0.49%  '
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('D48').Value = '0.0517'
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('D49').Value = '7.63'
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('D50').Value = '0.0967'
$ws.Range('E50').Value = '  +0.59%  '
$ws.Range('E51').Value = '  +0.31%  '

Write-Host "Applied cryptos update"
